{"js": "const body = context.document.body;\n\n// Grab all existing (original) paragraphs before we insert anything.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst lastOriginal = body.paragraphs.items[body.paragraphs.items.length - 1];\n\n// New OOXML content (the \"Birds of the Prairie\" / Wood Duck feature) that\n// replaces the entire old essay body.\nconst newContentOoxml = \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t xml:space=\\\"preserve\\\">The </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\\\"preserve\\\">Wood Duck </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">is one of the most stunningly pretty of all waterfowl. Males are iridescent chestnut and green, with ornate patterns on nearly every feather; the elegant females have a distinctive profile and delicate white pattern around the eye. These birds live in wooded swamps, where they nest in holes in trees or in nest boxes put up around lake margins. They are one of the few duck species equipped with strong claws that can grip bark and perch on branches. </w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Song</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t>Although better known for their distinctive coloration, Wood Duck calls are easily identifiable as well. Females make loud &#8220;</w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>oo</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\">-eek, </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>oo</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\">-eek&#8221; sounds when disturbed and when taking flight. Male Wood Ducks have a thin, rising and falling </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>zeeting</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> whistle. While flying, the wings of the wood duck make a whistling or whirring sound (Description adapted from text created by the Cornell Lab of Ornithology; https://www.alllaboutbirds.org/guide/wood_duck)</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\\\"preserve\\\">To hear an audio recording of the Wood Duck, captured by Paul Driver, visit </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Xeno</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>-canto (https://www.xeno-canto.org/182950), or download the fully interactive PDF of this issue&#8217;s Birds of the Prairie feature by using the hyperlink at the bottom of the web page.</w:t></w:r></w:p><w:p><w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/><w:bookmarkEnd w:id=\\\"0\\\"/></w:p><w:p/><w:p/><w:p><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n// Insert the new paragraphs right after the end of the current last\n// paragraph. Anchoring the insert on a RANGE (rather than the Body itself)\n// keeps the new paragraphs completely separate from the existing content,\n// so nothing gets merged together.\nconst insertionPoint = lastOriginal.getRange(Word.RangeLocation.end);\ninsertionPoint.insertOoxml(newContentOoxml, Word.InsertLocation.after);\nawait context.sync();\n\n// Now remove every paragraph that existed before our insertion, leaving\n// only the freshly inserted \"Birds of the Prairie\" content behind.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst totalParagraphs = body.paragraphs.items.length;\nconst newParagraphCount = 7; // paragraphs contributed by newContentOoxml\nconst originalCount = totalParagraphs - newParagraphCount;\nfor (let i = 0; i < originalCount; i++) {\n  body.paragraphs.load(\"items\");\n  await context.sync();\n  body.paragraphs.items[0].delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the entire body content with the new \"Birds of the Prairie\" /\n# Wood Duck feature content, preserving the section properties (sectPr)\n# which live outside the Content range.\n$newContentOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">The </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">Wood Duck </w:t></w:r><w:r><w:t xml:space=\"preserve\">is one of the most stunningly pretty of all waterfowl. Males are iridescent chestnut and green, with ornate patterns on nearly every feather; the elegant females have a distinctive profile and delicate white pattern around the eye. These birds live in wooded swamps, where they nest in holes in trees or in nest boxes put up around lake margins. They are one of the few duck species equipped with strong claws that can grip bark and perch on branches. </w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Song</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>Although better known for their distinctive coloration, Wood Duck calls are easily identifiable as well. Females make loud &#8220;</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>oo</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">-eek, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>oo</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">-eek&#8221; sounds when disturbed and when taking flight. Male Wood Ducks have a thin, rising and falling </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>zeeting</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> whistle. While flying, the wings of the wood duck make a whistling or whirring sound (Description adapted from text created by the Cornell Lab of Ornithology; https://www.alllaboutbirds.org/guide/wood_duck)</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">To hear an audio recording of the Wood Duck, captured by Paul Driver, visit </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Xeno</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>-canto (https://www.xeno-canto.org/182950), or download the fully interactive PDF of this issue&#8217;s Birds of the Prairie feature by using the hyperlink at the bottom of the web page.</w:t></w:r></w:p><w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p><w:p/><w:p/><w:p><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$d.Content.InsertXML($newContentOoxml)\n"}
